$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Make the header row bold (creates the new bold font + cellXf s="2") ---
$ws.Range("A1:D1").Font.Bold = $true

# --- Append the new journal entries (rows 4-13) -----------------------------
# Shared strings must be created in this exact order so the new <si> entries
# line up the same way they did in the authored workbook: the "Durée" values
# used further down the table, then the two Sqlite remarks for row 5, then
# the "Remarques" column top to bottom.
$ws.Range("B7").Value  = "3h45min"
$ws.Range("B9").Value  = "1h30min"
$ws.Range("B6").Value  = "1h45min"
$ws.Range("C5").Value  = "Sqlite"
$ws.Range("D5").Value  = "Création de la base de données selon le MLD"
$ws.Range("D4").Value  = "Chapitre conception"
$ws.Range("D6").Value  = "Création d'un formulaire de Login"
$ws.Range("D7").Value  = "Fin du login + création de compte fonctionnel"
$ws.Range("D8").Value  = "Login fonctionnel + vue du frigo"
$ws.Range("D9").Value  = "Ajout d'un ingrédient (Design formulaire) + upload d'image fonctionnel"
$ws.Range("D10").Value = "Ajout d'un ingrédient terminé"
$ws.Range("D11").Value = "Modification d'un ingrédient"
$ws.Range("D12").Value = "Modification d'un ingrédient terminé + événement clic sur une carte"
$ws.Range("D13").Value = "Suppression d'un aliment terminé + message date de péremption"

# Remaining cells (dates + already-existing shared strings) --------------
$ws.Range("A4").Value = 43445
$ws.Range("A5").Value = 43445
$ws.Range("A6").Value = 43445
$ws.Range("A7").Value = 43452
$ws.Range("A8").Value = 43473
$ws.Range("A9").Value = 43476
$ws.Range("A10").Value = 43480
$ws.Range("A11").Value = 43481
$ws.Range("A12").Value = 43483
$ws.Range("A13").Value = 43487

# Give the new date cells the same date formatting as the existing ones
$ws.Range("A2").Copy()
$ws.Range("A4:A13").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("B4").Value = "1h"
$ws.Range("C4").Value = "Documentation"

$ws.Range("B5").Value = "1h"

$ws.Range("C6").Value = "Code"

$ws.Range("B8").Value = "3h45min"
$ws.Range("C7").Value = "Code"
$ws.Range("C8").Value = "Code"

$ws.Range("B10").Value = "3h45min"
$ws.Range("B11").Value = "1h30min"
$ws.Range("B12").Value = "3h45min"
$ws.Range("B13").Value = "3h45min"

$ws.Range("C9").Value = "Code"
$ws.Range("C10").Value = "Code"
$ws.Range("C11").Value = "Code"
$ws.Range("C12").Value = "Code"
$ws.Range("C13").Value = "Code"

# --- Page setup (portrait, paper size 9 = A4) --------------------------------
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# --- Final selection, matching the authored workbook ------------------------
$ws.Range("A14").Select() | Out-Null
